# Update the shared "Enlace" link for the "Acceso compartido" requirement
# row (B5) with the new repository URL, then leave the sheet scrolled /
# selected / zoomed the way the author left it after editing that cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "Enlace: https://github.com/Dusbchek/Evidencia3"

$ws.Range("B5").Select()
$excel.ActiveWindow.Zoom = 67
